# Update the Users_db sheet: replace old mock login table with new ID/Password/Type data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers ---
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "Type"

# --- Data rows ---
$ws.Range("A2").Value = 311369318
$ws.Range("B2").Value = 111111
$ws.Range("C2").Value = 1

$ws.Range("A3").Value = 123456798
$ws.Range("B3").Value = 222222
$ws.Range("C3").Value = 2

$ws.Range("A4").Value = 456789123
$ws.Range("B4").Value = 3333
$ws.Range("C4").Value = 3

$ws.Range("A5").Value = 321456987
$ws.Range("B5").Value = 44444
$ws.Range("C5").Value = 1

# --- Font: Calibri -> Arial (entire used range) ---
$ws.Cells.Font.Name = "Arial"

# --- Column widths (closest achievable to 9.875 / 9.125 given engine's 1/6-char rounding) ---
$ws.Columns.Item(1).ColumnWidth = 9
$ws.Columns.Item(2).ColumnWidth = 8.333333

# --- Header alignment: A1:B1 center/center with border; C1 center/center without border ---
$headerRange = $ws.Range("A1:B1")
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4108    # xlCenter
$headerRange.Borders.LineStyle = 1

$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("C1").VerticalAlignment = -4108
$ws.Range("C1").Borders.LineStyle = -4142  # xlLineStyleNone

# --- Selection ---
$ws.Range("E10").Select() | Out-Null
